$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 947.3333
$ws.Range("I2").Value = 517.5
$ws.Range("J2").Value = 1162.25
$ws.Range("K2").Value = 517.5
$ws.Range("L2").Value = 1162.25
$ws.Range("M2").Value = -404.5
$ws.Range("N2").Value = -1388.25
$ws.Range("H19").Value = 906.9231
$ws.Range("J19").Value = 648.7778
$ws.Range("L19").Value = 648.7778
$ws.Range("N19").Value = -998.7778
$ws.Range("H38").Value = 2374.7778
$ws.Range("I38").Value = 72.166664
$ws.Range("K38").Value = 216.499992
$ws.Range("M38").Value = 155.500008
$ws.Range("H43").Value = 5749.5
$ws.Range("I43").Value = 5500.0
$ws.Range("J43").Value = 5999.0
$ws.Range("K43").Value = 5500.0
$ws.Range("L43").Value = 5999.0
$ws.Range("M43").Value = -5431.0
$ws.Range("N43").Value = -6137.0
$ws.Range("H44").Value = 0.0
$ws.Range("J44").Value = 0.0
$ws.Range("L44").Value = 0.0
$ws.Range("N44").Value = ""
$ws.Range("H58").Value = 1746.3529
$ws.Range("J58").Value = 2395.2856
$ws.Range("L58").Value = 7185.8568
$ws.Range("N58").Value = -7485.8568
$ws.Range("H64").Value = 6416.5
$ws.Range("I64").Value = 4833.0
$ws.Range("K64").Value = 4833.0
$ws.Range("M64").Value = -4585.0
$ws.Range("H67").Value = 6416.5
$ws.Range("I67").Value = 4833.0
$ws.Range("K67").Value = 4833.0
$ws.Range("M67").Value = -3975.0
$ws.Range("H80").Value = 711.3333
$ws.Range("I80").Value = 539.3333
$ws.Range("J80").Value = 883.3333
$ws.Range("K80").Value = 1617.9999
$ws.Range("L80").Value = 2649.9999
$ws.Range("M80").Value = -619.9999
$ws.Range("N80").Value = -4645.9999
$ws.Range("H83").Value = 711.3333
$ws.Range("I83").Value = 539.3333
$ws.Range("J83").Value = 883.3333
$ws.Range("K83").Value = 4853.9997
$ws.Range("L83").Value = 7949.9997
$ws.Range("M83").Value = 138.0002999999997
$ws.Range("N83").Value = -17933.9997
$ws.Range("H92").Value = 96.22222
$ws.Range("I92").Value = 98.25
$ws.Range("J92").Value = 80.0
$ws.Range("K92").Value = 98.25
$ws.Range("L92").Value = 80.0
$ws.Range("M92").Value = 1149.75
$ws.Range("N92").Value = -2576.0
$ws.Range("H98").Value = 766.1667
$ws.Range("I98").Value = 792.375
$ws.Range("J98").Value = 713.75
$ws.Range("K98").Value = 792.375
$ws.Range("L98").Value = 713.75
$ws.Range("M98").Value = 705.625
$ws.Range("N98").Value = -3709.75
$ws.Range("H111").Value = 822.2
$ws.Range("I111").Value = 609.6667
$ws.Range("J111").Value = 1141.0
$ws.Range("K111").Value = 1829.0001
$ws.Range("L111").Value = 3423.0
$ws.Range("M111").Value = 1237.9999
$ws.Range("N111").Value = -9557.0
$ws.Range("H116").Value = 3971.0
$ws.Range("J116").Value = 3950.0
$ws.Range("L116").Value = 3950.0
$ws.Range("N116").Value = -10834.0
$ws.Range("H122").Value = 766.1667
$ws.Range("I122").Value = 792.375
$ws.Range("J122").Value = 713.75
$ws.Range("K122").Value = 2377.125
$ws.Range("L122").Value = 2141.25
$ws.Range("M122").Value = 72.875
$ws.Range("N122").Value = -7041.25
$ws.Range("H127").Value = 593.0
$ws.Range("I127").Value = 593.0
$ws.Range("K127").Value = 1779.0
$ws.Range("M127").Value = 3181.0
$ws.Range("H138").Value = 2712.0908
$ws.Range("J138").Value = 3805.5
$ws.Range("L138").Value = 11416.5
$ws.Range("N138").Value = -21696.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 80.5
$ws.Range("I5").Value = 70.75
$ws.Range("K5").Value = 70.75
$ws.Range("M5").Value = 41.25
$ws.Range("H9").Value = 22666.0
$ws.Range("J9").Value = 22666.0
$ws.Range("L9").Value = 22666.0
$ws.Range("N9").Value = -23006.0
$ws.Range("H20").Value = 22666.0
$ws.Range("J20").Value = 22666.0
$ws.Range("L20").Value = 22666.0
$ws.Range("N20").Value = -23206.0
$ws.Range("H63").Value = 6299.625
$ws.Range("I63").Value = 1733.1666
$ws.Range("K63").Value = 1733.1666
$ws.Range("M63").Value = -1047.1666
$ws.Range("H66").Value = 6299.625
$ws.Range("I66").Value = 1733.1666
$ws.Range("K66").Value = 8665.833
$ws.Range("M66").Value = -5233.833000000001
$ws.Range("H122").Value = 1144.7142
$ws.Range("I122").Value = 1144.7142
$ws.Range("K122").Value = 3434.1426
$ws.Range("M122").Value = -984.1425999999997
$ws.Range("H132").Value = 293.75
$ws.Range("I132").Value = 202.0
$ws.Range("J132").Value = 446.66666
$ws.Range("K132").Value = 606.0
$ws.Range("L132").Value = 1339.99998
$ws.Range("M132").Value = 1924.0
$ws.Range("N132").Value = -6399.999980000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 80.5
$ws.Range("I4").Value = 70.75
$ws.Range("K4").Value = 70.75
$ws.Range("M4").Value = 44.25
$ws.Range("H95").Value = 5657.4287
$ws.Range("J95").Value = 5657.4287
$ws.Range("L95").Value = 5657.4287
$ws.Range("N95").Value = -11149.4287
$ws.Range("H99").Value = 4002.3333
$ws.Range("J99").Value = 3998.5
$ws.Range("L99").Value = 3998.5
$ws.Range("N99").Value = -6994.5
$ws.Range("H107").Value = 6259.1816
$ws.Range("I107").Value = 4856.375
$ws.Range("K107").Value = 4856.375
$ws.Range("M107").Value = -2936.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8078.0
$ws.Range("I4").Value = 13330.0
$ws.Range("J4").Value = 200.0
$ws.Range("K4").Value = 13330.0
$ws.Range("L4").Value = 200.0
$ws.Range("M4").Value = -13218.0
$ws.Range("N4").Value = -424.0
$ws.Range("H22").Value = 1933.9166
$ws.Range("I22").Value = 1642.2
$ws.Range("J22").Value = 2142.2856
$ws.Range("K22").Value = 1642.2
$ws.Range("L22").Value = 2142.2856
$ws.Range("M22").Value = -1292.2
$ws.Range("N22").Value = -2842.2856
$ws.Range("H105").Value = 1705.0
$ws.Range("I105").Value = 1705.0
$ws.Range("K105").Value = 1705.0
$ws.Range("M105").Value = 42.0
$ws.Range("H121").Value = 0.0
$ws.Range("J121").Value = 0.0
$ws.Range("L121").Value = 0.0
$ws.Range("N121").Value = ""
$ws.Range("H132").Value = 1891.25
$ws.Range("I132").Value = 1891.25
$ws.Range("K132").Value = 5673.75
$ws.Range("M132").Value = -3143.75
$ws.Range("H134").Value = 2981.0
$ws.Range("I134").Value = 2981.0
$ws.Range("K134").Value = 8943.0
$ws.Range("M134").Value = -6408.0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.6875
$ws.Range("I2").Value = 17.0
$ws.Range("K2").Value = 102.0
$ws.Range("M2").Value = 11.0
$ws.Range("H4").Value = 211313.9
$ws.Range("J4").Value = 22177.6
$ws.Range("L4").Value = 66532.79999999999
$ws.Range("N4").Value = -66756.79999999999
$ws.Range("H80").Value = 3496.0
$ws.Range("I80").Value = 3149.6
$ws.Range("J80").Value = 3743.4285
$ws.Range("K80").Value = 9448.8
$ws.Range("L80").Value = 11230.2855
$ws.Range("M80").Value = -8512.8
$ws.Range("N80").Value = -13102.2855
$ws.Range("H83").Value = 3496.0
$ws.Range("I83").Value = 3149.6
$ws.Range("J83").Value = 3743.4285
$ws.Range("K83").Value = 28346.4
$ws.Range("L83").Value = 33690.8565
$ws.Range("M83").Value = -23666.4
$ws.Range("N83").Value = -43050.8565
$ws.Range("H107").Value = 1362.75
$ws.Range("J107").Value = 1500.2858
$ws.Range("L107").Value = 4500.857400000001
$ws.Range("N107").Value = -8340.8574
$ws.Range("H141").Value = 15000.0
$ws.Range("I141").Value = 0.0
$ws.Range("K141").Value = 0.0
$ws.Range("M141").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 14472.333
$ws.Range("H46").Value = 1749.5
$ws.Range("I46").Value = 1749.5
$ws.Range("K46").Value = 1749.5
$ws.Range("M46").Value = -1593.5
$ws.Range("H57").Value = 26870.0
$ws.Range("J57").Value = 39792.0
$ws.Range("L57").Value = 39792.0
$ws.Range("N57").Value = -41432.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 92872.336
$ws.Range("J2").Value = 119019.0
$ws.Range("L2").Value = 119019.0
$ws.Range("N2").Value = -119243.0
$ws.Range("H40").Value = 4313.4443
$ws.Range("I40").Value = 4470.1665
$ws.Range("J40").Value = 4000.0
$ws.Range("K40").Value = 4470.1665
$ws.Range("L40").Value = 4000.0
$ws.Range("M40").Value = -4334.1665
$ws.Range("N40").Value = -4272.0
$ws.Range("H46").Value = 6222.8823
$ws.Range("I46").Value = 4474.75
$ws.Range("J46").Value = 6760.769
$ws.Range("K46").Value = 4474.75
$ws.Range("L46").Value = 6760.769
$ws.Range("M46").Value = -4286.75
$ws.Range("N46").Value = -7136.769
$ws.Range("H136").Value = 3831.6667
$ws.Range("J136").Value = 3831.6667
$ws.Range("L136").Value = 11495.0001
$ws.Range("N136").Value = -16595.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 479.33334
$ws.Range("I107").Value = 479.33334
$ws.Range("K107").Value = 1438.00002
$ws.Range("M107").Value = 481.9999800000001
$ws.Range("H126").Value = 4684.76
$ws.Range("I126").Value = 3136.4707
$ws.Range("J126").Value = 7974.875
$ws.Range("K126").Value = 9409.4121
$ws.Range("L126").Value = 23924.625
$ws.Range("M126").Value = -6939.4121
